# Auto-generated: apply scheduled-runner price refresh to Moogle_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 31250418
$ws.Range("I9").Value = 268.29166
$ws.Range("J9").Value = 125000870
$ws.Range("K9").Value = 268.29166
$ws.Range("L9").Value = 125000870
$ws.Range("M9").Value = -99.29165999999998
$ws.Range("N9").Value = -125001208
$ws.Range("H17").Value = 3574.3076
$ws.Range("J17").Value = 3574.3076
$ws.Range("L17").Value = 10722.9228
$ws.Range("N17").Value = -11058.9228
$ws.Range("H38").Value = 206.6
$ws.Range("I38").Value = 206.6
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 619.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -247.8
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 673.8570999999999
$ws.Range("I39").Value = 40.666668
$ws.Range("J39").Value = 1148.75
$ws.Range("K39").Value = 122.000004
$ws.Range("L39").Value = 3446.25
$ws.Range("M39").Value = 173.999996
$ws.Range("N39").Value = -4038.25
$ws.Range("H51").Value = 19666.666
$ws.Range("I51").Value = 19666.666
$ws.Range("K51").Value = 19666.666
$ws.Range("M51").Value = -19182.666
$ws.Range("H58").Value = 1116.1428
$ws.Range("I58").Value = 562.6
$ws.Range("K58").Value = 1687.8
$ws.Range("M58").Value = -1537.8
$ws.Range("H61").Value = 1998.3334
$ws.Range("I61").Value = 1998.3334
$ws.Range("K61").Value = 5995.0002
$ws.Range("M61").Value = -5823.0002
$ws.Range("H64").Value = 14999
$ws.Range("J64").Value = 14999
$ws.Range("L64").Value = 14999
$ws.Range("N64").Value = -15495
$ws.Range("H67").Value = 14999
$ws.Range("J67").Value = 14999
$ws.Range("L67").Value = 14999
$ws.Range("N67").Value = -16715
$ws.Range("H69").Value = 13757.451
$ws.Range("J69").Value = 13944.186
$ws.Range("L69").Value = 41832.558
$ws.Range("N69").Value = -43580.558
$ws.Range("H70").Value = 6103.3076
$ws.Range("I70").Value = 3474.1667
$ws.Range("K70").Value = 10422.5001
$ws.Range("M70").Value = -10152.5001
$ws.Range("H72").Value = 13757.451
$ws.Range("J72").Value = 13944.186
$ws.Range("L72").Value = 125497.674
$ws.Range("N72").Value = -134233.674
$ws.Range("H73").Value = 6103.3076
$ws.Range("I73").Value = 3474.1667
$ws.Range("K73").Value = 10422.5001
$ws.Range("M73").Value = -9486.500100000001
$ws.Range("H86").Value = 7360.5454
$ws.Range("I86").Value = 4496.7144
$ws.Range("K86").Value = 4496.7144
$ws.Range("M86").Value = -3373.7144
$ws.Range("H89").Value = 7360.5454
$ws.Range("I89").Value = 4496.7144
$ws.Range("K89").Value = 22483.572
$ws.Range("M89").Value = -16867.572
$ws.Range("H97").Value = 4579.1665
$ws.Range("J97").Value = 5095
$ws.Range("L97").Value = 15285
$ws.Range("N97").Value = -16277
$ws.Range("H98").Value = 1851.5385
$ws.Range("I98").Value = 1648.0476
$ws.Range("K98").Value = 1648.0476
$ws.Range("M98").Value = -150.0476000000001
$ws.Range("H113").Value = 10004120
$ws.Range("I113").Value = 50000000
$ws.Range("J113").Value = 5149.5
$ws.Range("K113").Value = 50000000
$ws.Range("L113").Value = 5149.5
$ws.Range("M113").Value = -49996746
$ws.Range("N113").Value = -11657.5
$ws.Range("H122").Value = 1851.5385
$ws.Range("I122").Value = 1648.0476
$ws.Range("K122").Value = 4944.142800000001
$ws.Range("M122").Value = -2494.142800000001
$ws.Range("H132").Value = 2167.7104
$ws.Range("I132").Value = 2172.2432
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6516.7296
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3986.7296
$ws.Range("N132").Value = -11060
$ws.Range("H137").Value = 2190.04
$ws.Range("I137").Value = 1757.8572
$ws.Range("K137").Value = 5273.571599999999
$ws.Range("M137").Value = -2723.571599999999
$ws.Range("H138").Value = 3713.1365
$ws.Range("I138").Value = 3244.4644
$ws.Range("J138").Value = 4533.3125
$ws.Range("K138").Value = 9733.393199999999
$ws.Range("L138").Value = 13599.9375
$ws.Range("M138").Value = -4593.393199999999
$ws.Range("N138").Value = -23879.9375
$ws.Range("H141").Value = 2206.4075
$ws.Range("I141").Value = 1869.2609
$ws.Range("K141").Value = 5607.7827
$ws.Range("M141").Value = -427.7826999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 9137
$ws.Range("I28").Value = 9137
$ws.Range("K28").Value = 9137
$ws.Range("M28").Value = -8945
$ws.Range("H31").Value = 980
$ws.Range("I31").Value = 980
$ws.Range("K31").Value = 980
$ws.Range("M31").Value = -686
$ws.Range("H32").Value = 6466.14
$ws.Range("I32").Value = 5226.9224
$ws.Range("J32").Value = 17619.1
$ws.Range("K32").Value = 5226.9224
$ws.Range("L32").Value = 17619.1
$ws.Range("M32").Value = -4939.9224
$ws.Range("N32").Value = -18193.1
$ws.Range("H37").Value = 44250
$ws.Range("I37").Value = 44250
$ws.Range("K37").Value = 44250
$ws.Range("M37").Value = -43977
$ws.Range("H45").Value = 4169728
$ws.Range("I45").Value = 7145439
$ws.Range("J45").Value = 3732.6
$ws.Range("K45").Value = 7145439
$ws.Range("L45").Value = 3732.6
$ws.Range("M45").Value = -7145062
$ws.Range("N45").Value = -4486.6
$ws.Range("H49").Value = 22916.666
$ws.Range("H56").Value = 9666.666999999999
$ws.Range("I56").Value = 7000
$ws.Range("J56").Value = 15000
$ws.Range("K56").Value = 7000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -6258
$ws.Range("N56").Value = -16484
$ws.Range("H61").Value = 8824.559999999999
$ws.Range("I61").Value = 7844.9473
$ws.Range("J61").Value = 11926.667
$ws.Range("K61").Value = 7844.9473
$ws.Range("L61").Value = 11926.667
$ws.Range("M61").Value = -7632.9473
$ws.Range("N61").Value = -12350.667
$ws.Range("H74").Value = 6486.5
$ws.Range("I74").Value = 3338.2144
$ws.Range("K74").Value = 3338.2144
$ws.Range("M74").Value = -2464.2144
$ws.Range("H77").Value = 6486.5
$ws.Range("I77").Value = 3338.2144
$ws.Range("K77").Value = 16691.072
$ws.Range("M77").Value = -12323.072
$ws.Range("H97").Value = 1703.1052
$ws.Range("I97").Value = 1647.0714
$ws.Range("K97").Value = 1647.0714
$ws.Range("M97").Value = -1151.0714
$ws.Range("H99").Value = 9137
$ws.Range("I99").Value = 9137
$ws.Range("K99").Value = 9137
$ws.Range("M99").Value = -6142
$ws.Range("H102").Value = 507
$ws.Range("I102").Value = 507
$ws.Range("K102").Value = 507
$ws.Range("M102").Value = 1115
$ws.Range("H110").Value = 2052.8462
$ws.Range("I110").Value = 1807.25
$ws.Range("K110").Value = 1807.25
$ws.Range("M110").Value = 237.75
$ws.Range("H120").Value = 60777
$ws.Range("J120").Value = 60777
$ws.Range("L120").Value = 60777
$ws.Range("N120").Value = -70453
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 37589.5
$ws.Range("J130").Value = 37589.5
$ws.Range("L130").Value = 37589.5
$ws.Range("N130").Value = -47629.5
$ws.Range("H132").Value = 4469.593
$ws.Range("I132").Value = 2303.5908
$ws.Range("K132").Value = 6910.7724
$ws.Range("M132").Value = -4380.7724
$ws.Range("H136").Value = 8824.559999999999
$ws.Range("I136").Value = 7844.9473
$ws.Range("J136").Value = 11926.667
$ws.Range("K136").Value = 23534.8419
$ws.Range("L136").Value = 35780.001
$ws.Range("M136").Value = -20984.8419
$ws.Range("N136").Value = -40880.001
$ws.Range("H139").Value = 121385.336
$ws.Range("J139").Value = 121385.336
$ws.Range("L139").Value = 121385.336
$ws.Range("N139").Value = -131665.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 10175
$ws.Range("I75").Value = 10175
$ws.Range("K75").Value = 10175
$ws.Range("M75").Value = -9239
$ws.Range("H78").Value = 10175
$ws.Range("I78").Value = 10175
$ws.Range("K78").Value = 30525
$ws.Range("M78").Value = -25845
$ws.Range("H102").Value = 9894.333000000001
$ws.Range("I102").Value = 9894.333000000001
$ws.Range("K102").Value = 9894.333000000001
$ws.Range("M102").Value = -6649.333000000001
$ws.Range("H105").Value = 731239.1
$ws.Range("I105").Value = 2863471.2
$ws.Range("J105").Value = 3159.9268
$ws.Range("K105").Value = 2863471.2
$ws.Range("L105").Value = 3159.9268
$ws.Range("M105").Value = -2861724.2
$ws.Range("N105").Value = -6653.9268
$ws.Range("H107").Value = 4588.875
$ws.Range("J107").Value = 5999.5
$ws.Range("L107").Value = 5999.5
$ws.Range("N107").Value = -9839.5
$ws.Range("H134").Value = 4336.2085
$ws.Range("I134").Value = 3160.6667
$ws.Range("J134").Value = 6922.4
$ws.Range("K134").Value = 9482.000100000001
$ws.Range("L134").Value = 20767.2
$ws.Range("M134").Value = -6947.000100000001
$ws.Range("N134").Value = -25837.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5817.6
$ws.Range("J29").Value = 7267.25
$ws.Range("L29").Value = 7267.25
$ws.Range("N29").Value = -7853.25
$ws.Range("H31").Value = 6412.896
$ws.Range("I31").Value = 2726.818
$ws.Range("K31").Value = 2726.818
$ws.Range("M31").Value = -2431.818
$ws.Range("H34").Value = 6412.896
$ws.Range("I34").Value = 2726.818
$ws.Range("K34").Value = 2726.818
$ws.Range("M34").Value = -2524.818
$ws.Range("H51").Value = 54995
$ws.Range("J51").Value = 54995
$ws.Range("L51").Value = 54995
$ws.Range("N51").Value = -56467
$ws.Range("H58").Value = 6106.7334
$ws.Range("I58").Value = 5682.647
$ws.Range("J58").Value = 6661.3076
$ws.Range("K58").Value = 5682.647
$ws.Range("L58").Value = 6661.3076
$ws.Range("M58").Value = -5479.647
$ws.Range("N58").Value = -7067.3076
$ws.Range("H61").Value = 54995
$ws.Range("J61").Value = 54995
$ws.Range("L61").Value = 54995
$ws.Range("N61").Value = -55691
$ws.Range("H86").Value = 2861331.5
$ws.Range("I86").Value = 4003880.5
$ws.Range("J86").Value = 4958.75
$ws.Range("K86").Value = 4003880.5
$ws.Range("L86").Value = 4958.75
$ws.Range("M86").Value = -4002757.5
$ws.Range("N86").Value = -7204.75
$ws.Range("H89").Value = 2861331.5
$ws.Range("I89").Value = 4003880.5
$ws.Range("J89").Value = 4958.75
$ws.Range("K89").Value = 20019402.5
$ws.Range("L89").Value = 24793.75
$ws.Range("M89").Value = -20013786.5
$ws.Range("N89").Value = -36025.75
$ws.Range("H95").Value = 47965.832
$ws.Range("J95").Value = 47965.832
$ws.Range("L95").Value = 47965.832
$ws.Range("N95").Value = -53457.832
$ws.Range("H99").Value = 8633.817999999999
$ws.Range("J99").Value = 5975.125
$ws.Range("L99").Value = 5975.125
$ws.Range("N99").Value = -8971.125
$ws.Range("H122").Value = 3116.0417
$ws.Range("J122").Value = 3820.3333
$ws.Range("L122").Value = 11460.9999
$ws.Range("N122").Value = -16360.9999
$ws.Range("H126").Value = 8633.817999999999
$ws.Range("J126").Value = 5975.125
$ws.Range("L126").Value = 17925.375
$ws.Range("N126").Value = -22865.375
$ws.Range("H132").Value = 8126.625
$ws.Range("I132").Value = 6574.15
$ws.Range("K132").Value = 19722.45
$ws.Range("M132").Value = -17192.45
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 6106.7334
$ws.Range("I136").Value = 5682.647
$ws.Range("J136").Value = 6661.3076
$ws.Range("K136").Value = 17047.941
$ws.Range("L136").Value = 19983.9228
$ws.Range("M136").Value = -14497.941
$ws.Range("N136").Value = -25083.9228
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 92499.5
$ws.Range("J138").Value = 92499.5
$ws.Range("L138").Value = 92499.5
$ws.Range("N138").Value = -102779.5
$ws.Range("H140").Value = 83970.39999999999
$ws.Range("J140").Value = 83970.39999999999
$ws.Range("L140").Value = 83970.39999999999
$ws.Range("N140").Value = -94330.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 49062.61
$ws.Range("I11").Value = 75191.336
$ws.Range("K11").Value = 225574.008
$ws.Range("M11").Value = -225434.008
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 2666.3333
$ws.Range("J32").Value = 3052.6843
$ws.Range("K32").Value = 7998.999899999999
$ws.Range("L32").Value = 9158.052899999999
$ws.Range("M32").Value = -7715.999899999999
$ws.Range("N32").Value = -9724.052899999999
$ws.Range("H39").Value = 19499
$ws.Range("J39").Value = 19499
$ws.Range("L39").Value = 58497
$ws.Range("N39").Value = -59085
$ws.Range("H55").Value = 3249.5
$ws.Range("J55").Value = 4999.5
$ws.Range("L55").Value = 14998.5
$ws.Range("N55").Value = -15352.5
$ws.Range("H57").Value = 14999.8
$ws.Range("I57").Value = 7499.5
$ws.Range("K57").Value = 22498.5
$ws.Range("M57").Value = -21939.5
$ws.Range("H68").Value = 6484.5713
$ws.Range("J68").Value = 1400
$ws.Range("L68").Value = 4200
$ws.Range("N68").Value = -5822
$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 15000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -16622
$ws.Range("H70").Value = 10157.9
$ws.Range("I70").Value = 3596.5
$ws.Range("K70").Value = 10789.5
$ws.Range("M70").Value = -10474.5
$ws.Range("H71").Value = 6484.5713
$ws.Range("J71").Value = 1400
$ws.Range("L71").Value = 12600
$ws.Range("N71").Value = -20712
$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 45000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -53112
$ws.Range("H73").Value = 10157.9
$ws.Range("I73").Value = 3596.5
$ws.Range("K73").Value = 10789.5
$ws.Range("M73").Value = -9697.5
$ws.Range("H75").Value = 10941.083
$ws.Range("J75").Value = 11390.272
$ws.Range("L75").Value = 34170.81600000001
$ws.Range("N75").Value = -36166.81600000001
$ws.Range("H78").Value = 10941.083
$ws.Range("J78").Value = 11390.272
$ws.Range("L78").Value = 102512.448
$ws.Range("N78").Value = -112496.448
$ws.Range("H80").Value = 2286.5
$ws.Range("J80").Value = 2417.1428
$ws.Range("L80").Value = 7251.428400000001
$ws.Range("N80").Value = -9123.428400000001
$ws.Range("H83").Value = 2286.5
$ws.Range("J83").Value = 2417.1428
$ws.Range("L83").Value = 21754.2852
$ws.Range("N83").Value = -31114.2852
$ws.Range("H131").Value = 6184.25
$ws.Range("J131").Value = 7141.933
$ws.Range("L131").Value = 21425.799
$ws.Range("N131").Value = -31505.799
$ws.Range("H132").Value = 2887.3125
$ws.Range("I132").Value = 2974.875
$ws.Range("J132").Value = 2799.75
$ws.Range("K132").Value = 26773.875
$ws.Range("L132").Value = 25197.75
$ws.Range("M132").Value = -24243.875
$ws.Range("N132").Value = -30257.75
$ws.Range("H133").Value = 7787.2666
$ws.Range("I133").Value = 6423.222
$ws.Range("J133").Value = 9833.333000000001
$ws.Range("K133").Value = 19269.666
$ws.Range("L133").Value = 29499.999
$ws.Range("M133").Value = -14209.666
$ws.Range("N133").Value = -39619.999
$ws.Range("H134").Value = 4176.2856
$ws.Range("I134").Value = 3036
$ws.Range("J134").Value = 19000
$ws.Range("K134").Value = 9108
$ws.Range("L134").Value = 57000
$ws.Range("M134").Value = -4038
$ws.Range("N134").Value = -67140
$ws.Range("H136").Value = 4035
$ws.Range("I136").Value = 3716
$ws.Range("J136").Value = 4513.5
$ws.Range("K136").Value = 11148
$ws.Range("L136").Value = 13540.5
$ws.Range("M136").Value = -6048
$ws.Range("N136").Value = -23740.5
$ws.Range("H137").Value = 3438.2666
$ws.Range("J137").Value = 4533.25
$ws.Range("L137").Value = 13599.75
$ws.Range("N137").Value = -23799.75
$ws.Range("H139").Value = 4675.731
$ws.Range("J139").Value = 6245
$ws.Range("L139").Value = 18735
$ws.Range("N139").Value = -29015
$ws.Range("H141").Value = 4636.6665
$ws.Range("I141").Value = 4240.3335
$ws.Range("J141").Value = 5033
$ws.Range("K141").Value = 12721.0005
$ws.Range("L141").Value = 15099
$ws.Range("M141").Value = -7541.000499999998
$ws.Range("N141").Value = -25459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4922.154
$ws.Range("I70").Value = 4000
$ws.Range("K70").Value = 4000
$ws.Range("M70").Value = -3730
$ws.Range("H73").Value = 4922.154
$ws.Range("I73").Value = 4000
$ws.Range("K73").Value = 4000
$ws.Range("M73").Value = -3064
$ws.Range("H97").Value = 762.2
$ws.Range("J97").Value = 647
$ws.Range("L97").Value = 647
$ws.Range("N97").Value = -1639
$ws.Range("H102").Value = 5403.077
$ws.Range("I102").Value = 3229
$ws.Range("J102").Value = 12650
$ws.Range("K102").Value = 3229
$ws.Range("L102").Value = 12650
$ws.Range("M102").Value = -1607
$ws.Range("N102").Value = -15894
$ws.Range("H132").Value = 5488.341
$ws.Range("I132").Value = 4733.4546
$ws.Range("J132").Value = 7753
$ws.Range("K132").Value = 14200.3638
$ws.Range("L132").Value = 23259
$ws.Range("M132").Value = -11670.3638
$ws.Range("N132").Value = -28319

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1223.6875
$ws.Range("I16").Value = 1038.6
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1038.6
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -868.5999999999999
$ws.Range("N16").Value = -4340
$ws.Range("H40").Value = 8002.1
$ws.Range("I40").Value = 4578
$ws.Range("K40").Value = 4578
$ws.Range("M40").Value = -4442
$ws.Range("H55").Value = 984.2
$ws.Range("I55").Value = 440.66666
$ws.Range("K55").Value = 440.66666
$ws.Range("M55").Value = -267.66666
$ws.Range("H68").Value = 10637.218
$ws.Range("I68").Value = 8785.666999999999
$ws.Range("J68").Value = 11827.5
$ws.Range("K68").Value = 8785.666999999999
$ws.Range("L68").Value = 11827.5
$ws.Range("M68").Value = -8036.666999999999
$ws.Range("N68").Value = -13325.5
$ws.Range("H71").Value = 10637.218
$ws.Range("I71").Value = 8785.666999999999
$ws.Range("J71").Value = 11827.5
$ws.Range("K71").Value = 43928.335
$ws.Range("L71").Value = 59137.5
$ws.Range("M71").Value = -40184.335
$ws.Range("N71").Value = -66625.5
$ws.Range("H82").Value = 2379.6
$ws.Range("I82").Value = 1876.3334
$ws.Range("K82").Value = 1876.3334
$ws.Range("M82").Value = -1515.3334
$ws.Range("H85").Value = 2379.6
$ws.Range("I85").Value = 1876.3334
$ws.Range("K85").Value = 1876.3334
$ws.Range("M85").Value = -628.3334
$ws.Range("H100").Value = 4842.2104
$ws.Range("J100").Value = 8502
$ws.Range("L100").Value = 8502
$ws.Range("N100").Value = -9584
$ws.Range("H122").Value = 3993.6667
$ws.Range("I122").Value = 2057.8
$ws.Range("J122").Value = 7220.1113
$ws.Range("K122").Value = 6173.400000000001
$ws.Range("L122").Value = 21660.3339
$ws.Range("M122").Value = -3723.400000000001
$ws.Range("N122").Value = -26560.3339
$ws.Range("H132").Value = 8541.346
$ws.Range("I132").Value = 9159.857
$ws.Range("K132").Value = 27479.571
$ws.Range("M132").Value = -24949.571
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 8988.151
$ws.Range("I136").Value = 4806.143
$ws.Range("K136").Value = 14418.429
$ws.Range("M136").Value = -11868.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 46986.43
$ws.Range("J54").Value = 45781.4
$ws.Range("L54").Value = 45781.4
$ws.Range("N54").Value = -46821.4
$ws.Range("H62").Value = 9248.75
$ws.Range("I62").Value = 4000
$ws.Range("K62").Value = 4000
$ws.Range("M62").Value = -3376
$ws.Range("H65").Value = 9248.75
$ws.Range("I65").Value = 4000
$ws.Range("K65").Value = 20000
$ws.Range("M65").Value = -16880
$ws.Range("H75").Value = 87016.86
$ws.Range("I75").Value = 84779.5
$ws.Range("K75").Value = 84779.5
$ws.Range("M75").Value = -83843.5
$ws.Range("H78").Value = 87016.86
$ws.Range("I78").Value = 84779.5
$ws.Range("K78").Value = 254338.5
$ws.Range("M78").Value = -249658.5
$ws.Range("H126").Value = 4316.6665
$ws.Range("I126").Value = 1475
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 4425
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
$ws.Range("M126").Value = -1955
$ws.Range("H132").Value = 3274.8276
$ws.Range("I132").Value = 2287.3076
$ws.Range("K132").Value = 6861.9228
$ws.Range("M132").Value = -4331.9228
$ws.Range("H136").Value = 12012.625
$ws.Range("I136").Value = 22055.5
$ws.Range("K136").Value = 66166.5
$ws.Range("M136").Value = -63616.5
